# Generate Report for Handoff
# Update status from "Ready for handoff" to "In Translation" and refresh
# the associated handoff/generate timestamps across the Overview, zh-cn,
# and de-de sheets. Also shrink the now-shorter "Status" column widths.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-13 12:45:24"

# Status column shrank (shorter text), so the "Status"-related columns
# are narrowed to match (target XML width ~13.41 chars).
$newStatusColWidth = 13.4101848602295 - 5/6
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2).
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-13 12:45:14"
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2).
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-13 12:45:24"
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
